$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Source" text (column D) to previously-blank cells in rows 57-70 ---
$ws.Range("D57").Value = "Developed by RAPID Team"
$ws.Range("D58").Value = "Developed by RAPID Team"
$ws.Range("D59").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D60").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D61").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D62").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D63").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D64").Value = "Developed by RAPID Team"
$ws.Range("D65").Value = "Developed by RAPID Team"

# Row 66 gets a brand-new source string not previously used anywhere in the
# workbook, and its row grows taller to accommodate the wrapped text.
$ws.Range("D66").Value = "RAPID Team Modified from National Compensation Survey"
$ws.Rows.Item(66).RowHeight = 96

$ws.Range("D67").Value = "Developed by RAPID Team"
$ws.Range("D68").Value = "RAPID Team Modified"
$ws.Range("D69").Value = "Developed by RAPID Team"
$ws.Range("D70").Value = "Developed by RAPID Team"

# --- Update the saved view state (scroll position / zoom / selection) ---
$win = $excel.ActiveWindow
$win.Zoom = 125
[void]$ws.Range("D69").Select()
$win.ScrollRow = 92
$win.ScrollColumn = 1
